$d = $word.ActiveDocument

$replacements = @(
    @{old = "871÷3=290, 1"; new = "396÷4=99, 0"},
    @{old = "773÷6=128, 5"; new = "944÷3=314, 2"},
    @{old = "693÷3=231, 0"; new = "930÷9=103, 3"},
    @{old = "213÷2=106, 1"; new = "147÷3=49, 0"},
    @{old = "358÷8=44, 6"; new = "962÷6=160, 2"},
    @{old = "973÷9=108, 1"; new = "671÷8=83, 7"},
    @{old = "293÷8=36, 5"; new = "906÷3=302, 0"},
    @{old = "828÷2=414, 0"; new = "555÷5=111, 0"},
    @{old = "386÷5=77, 1"; new = "316÷2=158, 0"},
    @{old = "403÷5=80, 3"; new = "556÷6=92, 4"},
    @{old = "132÷6=22, 0"; new = "185÷6=30, 5"},
    @{old = "203÷9=22, 5"; new = "504÷8=63, 0"},
    @{old = "986÷7=140, 6"; new = "375÷2=187, 1"},
    @{old = "445÷6=74, 1"; new = "422÷5=84, 2"},
    @{old = "841÷3=280, 1"; new = "526÷7=75, 1"},
    @{old = "935÷3=311, 2"; new = "254÷2=127, 0"},
    @{old = "235÷2=117, 1"; new = "490÷4=122, 2"},
    @{old = "208÷5=41, 3"; new = "514÷2=257, 0"},
    @{old = "689÷6=114, 5"; new = "199÷2=99, 1"},
    @{old = "598÷7=85, 3"; new = "200÷2=100, 0"},
    @{old = "899÷2=449, 1"; new = "531÷5=106, 1"},
    @{old = "604÷4=151, 0"; new = "152÷7=21, 5"},
    @{old = "711÷8=88, 7"; new = "670÷8=83, 6"},
    @{old = "154÷7=22, 0"; new = "545÷4=136, 1"},
    @{old = "424÷4=106, 0"; new = "723÷2=361, 1"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
